$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet and name it "nr_studies"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "nr_studies"

# Header row
$ws.Range("A1").Value = "outcome"
$ws.Range("B1").Value = "n_effect_sizes"
$ws.Range("C1").Value = "k_studies"

# Data rows
$ws.Range("A2").Value = "NS"
$ws.Range("B2").Value = 722
$ws.Range("C2").Value = 81

$ws.Range("A3").Value = "NT"
$ws.Range("B3").Value = 381
$ws.Range("C3").Value = 51

# Style header row: bold + centered (matches the other sheets in the workbook)
$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108

# Restore original active sheet/selection (adding a sheet makes it active by default)
$wb.Worksheets.Item(1).Activate()
$wb.Worksheets.Item(1).Range("A1").Select() | Out-Null
